$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "as of" header date from June 1 to June 2
$ws.Range("E2").Value = "（6月2日時点）"

# 2) Insert a new row for 2021-06-02 right below the grand-total row (row 4),
#    pushing the existing data rows (and everything below) down by one.
$ws.Rows.Item(5).Insert()

# Copy the formatting (styles only) from the row directly below (the row that
# used to be row 5) onto the freshly inserted row so it matches the rest of
# the data rows exactly.
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Range("A5").Value = 44349
$ws.Range("B5").Value = "(水)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 71456
$ws.Range("E5").Value = 87670

# 3) Update the running grand-total row (row 4) to include the new day's counts.
$ws.Range("D4").Value = 4725022
$ws.Range("E4").Value = 3227298
